$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the confidential disclosure text with the new "as of" date
$ws.Range("A40").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-13 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for each holding row
$ws.Range("D2").Value = 0.03068841813263206
$ws.Range("E2").Value = 0.002380345150046814
$ws.Range("D3").Value = 0.03355410165629352
$ws.Range("E3").Value = 0.008013995334888335
$ws.Range("D4").Value = 0.0338448069688772
$ws.Range("E4").Value = -0.006056601695848429
$ws.Range("D5").Value = 0.06929501539154186
$ws.Range("E5").Value = -0.01107689825335145
$ws.Range("D6").Value = 0.03012154277309388
$ws.Range("E6").Value = -0.002672605790645877
$ws.Range("D7").Value = 0.01592953303222971
$ws.Range("E7").Value = -0.009054537797430928
$ws.Range("D8").Value = 0.03187676927259754
$ws.Range("E8").Value = 0.001882391455579668
$ws.Range("D9").Value = 0.03229083799346994
$ws.Range("E9").Value = -0.02060237417835786
$ws.Range("D10").Value = 0.05006374086677515
$ws.Range("E10").Value = 0.00627198451545663
$ws.Range("D11").Value = 0.02680377521845292
$ws.Range("E11").Value = 0.01961261436636175
$ws.Range("D12").Value = 0.01503635959427228
$ws.Range("E12").Value = -0.02992973019866396
$ws.Range("D13").Value = 0.01586580148293252
$ws.Range("E13").Value = 0.001761804087385466
$ws.Range("D14").Value = 0.01484982368536442
$ws.Range("E14").Value = -0.008432888264230498
$ws.Range("D15").Value = 0.007076997215378533
$ws.Range("E15").Value = -0.01050635911209408
$ws.Range("D16").Value = 0.007379442550201173
$ws.Range("E16").Value = -0.01737373737373737
$ws.Range("D17").Value = 0.03168632002294335
$ws.Range("E17").Value = 0.004569593676670536
$ws.Range("D18").Value = 0.02804188169076446
$ws.Range("E18").Value = -0.003728070175438858
$ws.Range("D19").Value = 0.03014576821580919
$ws.Range("E19").Value = -0.003548247511899638
$ws.Range("D20").Value = 0.03290560518976069
$ws.Range("E20").Value = -0.003029788197984162
$ws.Range("D21").Value = 0.04866108773355889
$ws.Range("E21").Value = 0.003297233540638933
$ws.Range("D22").Value = 0.0283624029328439
$ws.Range("E22").Value = -0.01187910643889611
$ws.Range("D23").Value = 0.02997805361239552
$ws.Range("E23").Value = -0.01457698762976323
$ws.Range("D24").Value = 0.02774838113479055
$ws.Range("E24").Value = 0.008824418253248867
$ws.Range("D25").Value = 0.0121559544554223
$ws.Range("E25").Value = -0.03715967623252392
$ws.Range("D26").Value = 0.01274742795679448
$ws.Range("E26").Value = -0.02157705464432957
$ws.Range("D27").Value = 0.02888641789373192
$ws.Range("E27").Value = -0.004541583877377242
$ws.Range("D28").Value = 0.02831916983507505
$ws.Range("E28").Value = -0.003369130343230275
$ws.Range("D29").Value = 0.03100707587911802
$ws.Range("E29").Value = 0.01319774989182165
$ws.Range("D30").Value = 0.03317059426315427
$ws.Range("E30").Value = -0.001421332344580417
$ws.Range("D31").Value = 0.03057996268909122
$ws.Range("E31").Value = -0.01279707495429616
$ws.Range("D32").Value = 0.0278611226181964
$ws.Range("E32").Value = 0.02307537957327277
$ws.Range("D33").Value = 0.02970766040178083
$ws.Range("E33").Value = 0.01035635651961186
$ws.Range("D34").Value = 0.03101154826854239
$ws.Range("E34").Value = 0.006009037592539235
$ws.Range("D35").Value = 0.02944323037706529
$ws.Range("E35").Value = -0.01487341772151896
$ws.Range("D36").Value = 0.0329033689950485
$ws.Range("E36").Value = -0.001942594354582994
$ws.Range("E37").Value = -0.002253338871642363

$ws.Protect()
